$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("F4").Value = 772
$ws.Range("F5").Value = 2505
$ws.Range("F6").Value = 59
$ws.Range("F7").Value = 1853
$ws.Range("F8").Value = 3195
$ws.Range("F9").Value = 199
$ws.Range("F10").Value = 4687
$ws.Range("F11").Value = 437
$ws.Range("F12").Value = 261
$ws.Range("F13").Value = 151
$ws.Range("F14").Value = 606
$ws.Range("F16").Value = 12
$ws.Range("F17").Value = 6
$ws.Range("F18").Value = 638
$ws.Range("F19").Value = 278
$ws.Range("F21").Value = 86
$ws.Range("F22").Value = 135
$ws.Range("F23").Value = 327
$ws.Range("F24").Value = 4659
$ws.Range("F26").Value = 30
$ws.Range("F28").Value = 5667
$ws.Range("F30").Value = 1172
$ws.Range("F32").Value = 644
$ws.Range("F33").Value = 4401
$ws.Range("F35").Value = 66
$ws.Range("F37").Value = 772
$ws.Range("F38").Value = 51
$ws.Range("F39").Value = 705
$ws.Range("F40").Value = 712

$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 33
$ws.Range("F4").Value = 9
$ws.Range("F5").Value = 8
$ws.Range("F6").Value = 43

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 216
$ws.Range("F4").Value = 31

$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 216
$ws.Range("F5").Value = 31
$ws.Range("F7").Value = 772
$ws.Range("F8").Value = 2505
$ws.Range("F9").Value = 59
$ws.Range("F10").Value = 1853
$ws.Range("F12").Value = 3195
$ws.Range("F13").Value = 199
$ws.Range("F14").Value = 4687
$ws.Range("F15").Value = 437
$ws.Range("F16").Value = 261
$ws.Range("F17").Value = 151
$ws.Range("F18").Value = 606
$ws.Range("F20").Value = 12
$ws.Range("F21").Value = 6
$ws.Range("F22").Value = 638
$ws.Range("F23").Value = 278
$ws.Range("F25").Value = 33
$ws.Range("F26").Value = 86
$ws.Range("F27").Value = 135
$ws.Range("F28").Value = 327
$ws.Range("F29").Value = 4659
$ws.Range("F31").Value = 30
$ws.Range("F33").Value = 5667
$ws.Range("F35").Value = 1172
$ws.Range("F37").Value = 644
$ws.Range("F38").Value = 4401
$ws.Range("F40").Value = 9
$ws.Range("F41").Value = 66
$ws.Range("F43").Value = 772
$ws.Range("F44").Value = 51
$ws.Range("F45").Value = 705
$ws.Range("F46").Value = 712
$ws.Range("F47").Value = 8
$ws.Range("F48").Value = 43
